# B6-PowerPoint.pptx - Apr 21 2020 commit
#
# 1) Three tables (slides 14-16) get a new built-in table style applied
#    (Table Design ribbon -> style gallery pick), changing their
#    <a:tableStyleId> from the local "Table_0" style
#    {CDD7A4D9-40D5-4805-A191-EDB4B8677A5A} to the built-in style
#    {1827DE7F-8DDE-46E9-99CF-AD639B7788DD}.
# 2) The deck's Design/Theme was reset to the default Office Theme from
#    the Design tab; PowerPoint keeps the presentation's previous theme
#    ("Integral") around for re-use by the Notes Master, which is why the
#    theme1.xml / theme2.xml parts end up swapped in the saved package.

$p = $ppt.ActivePresentation

$newTableStyleId = "{1827DE7F-8DDE-46E9-99CF-AD639B7788DD}"

for ($idx = 1; $idx -le $p.Slides.Count; $idx++) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# Re-apply the default Office Theme design to the whole deck (Design tab).
# PowerPoint hands the slide master the new "Office Theme" theme and
# relocates the deck's former theme ("Integral") onto the Notes Master.
$master = $p.Slides.Item(1).Master
$notesMaster = $p.NotesMaster
$slideTheme = $master.Theme
$notesTheme = $notesMaster.Theme

$notesMaster.Theme = $slideTheme
$master.Theme = $notesTheme
